$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7777777777777778
$ws.Range("C2").Value = 0.7142857142857143
$ws.Range("D2").Value = 0.7446808510638298
$ws.Range("E2").Value = 49

# Row 3
$ws.Range("B3").Value = 0.7083333333333334
$ws.Range("C3").Value = 0.7727272727272727
$ws.Range("D3").Value = 0.7391304347826088
$ws.Range("E3").Value = 44

# Row 4
$ws.Range("B4").Value = 0.7419354838709677
$ws.Range("C4").Value = 0.7419354838709677
$ws.Range("D4").Value = 0.7419354838709677
$ws.Range("E4").Value = 0.7419354838709677

# Row 5
$ws.Range("B5").Value = 0.7430555555555556
$ws.Range("C5").Value = 0.7435064935064934
$ws.Range("D5").Value = 0.7419056429232193
$ws.Range("E5").Value = 93

# Row 6
$ws.Range("B6").Value = 0.7449223416965354
$ws.Range("C6").Value = 0.7419354838709677
$ws.Range("D6").Value = 0.7420548476619618
$ws.Range("E6").Value = 93

# Row 7
$ws.Range("B7").Value = 0.6111111111111112
$ws.Range("C7").Value = 0.8979591836734694
$ws.Range("D7").Value = 0.7272727272727272
$ws.Range("E7").Value = 49

# Row 8
$ws.Range("B8").Value = 0.7619047619047619
$ws.Range("C8").Value = 0.3636363636363636
$ws.Range("D8").Value = 0.4923076923076923
$ws.Range("E8").Value = 44

# Row 9
$ws.Range("B9").Value = 0.6451612903225806
$ws.Range("C9").Value = 0.6451612903225806
$ws.Range("D9").Value = 0.6451612903225806
$ws.Range("E9").Value = 0.6451612903225806

# Row 10
$ws.Range("B10").Value = 0.6865079365079365
$ws.Range("C10").Value = 0.6307977736549165
$ws.Range("D10").Value = 0.6097902097902097
$ws.Range("E10").Value = 93

# Row 11
$ws.Range("B11").Value = 0.6824543437446663
$ws.Range("C11").Value = 0.6451612903225806
$ws.Range("D11").Value = 0.6161064741709903
$ws.Range("E11").Value = 93

# Row 12
$ws.Range("B12").Value = 0.6842105263157895
$ws.Range("C12").Value = 0.7959183673469388
$ws.Range("D12").Value = 0.7358490566037734
$ws.Range("E12").Value = 49

# Row 13
$ws.Range("B13").Value = 0.7222222222222222
$ws.Range("C13").Value = 0.5909090909090909
$ws.Range("D13").Value = 0.65
$ws.Range("E13").Value = 44

# Row 14
$ws.Range("B14").Value = 0.6989247311827957
$ws.Range("C14").Value = 0.6989247311827957
$ws.Range("D14").Value = 0.6989247311827957
$ws.Range("E14").Value = 0.6989247311827957

# Row 15
$ws.Range("B15").Value = 0.7032163742690059
$ws.Range("C15").Value = 0.6934137291280149
$ws.Range("D15").Value = 0.6929245283018868
$ws.Range("E15").Value = 93

# Row 16
$ws.Range("B16").Value = 0.702194554486575
$ws.Range("C16").Value = 0.6989247311827957
$ws.Range("D16").Value = 0.6952322986406978
$ws.Range("E16").Value = 93

# Row 17
$ws.Range("B17").Value = 0.6744186046511628
$ws.Range("C17").Value = 0.5918367346938775
$ws.Range("D17").Value = 0.6304347826086958
$ws.Range("E17").Value = 49

# Row 18
$ws.Range("B18").Value = 0.6
$ws.Range("C18").Value = 0.6818181818181818
$ws.Range("D18").Value = 0.6382978723404256
$ws.Range("E18").Value = 44

# Row 19
$ws.Range("B19").Value = 0.6344086021505376
$ws.Range("C19").Value = 0.6344086021505376
$ws.Range("D19").Value = 0.6344086021505376
$ws.Range("E19").Value = 0.6344086021505376

# Row 20
$ws.Range("B20").Value = 0.6372093023255814
$ws.Range("C20").Value = 0.6368274582560296
$ws.Range("D20").Value = 0.6343663274745607
$ws.Range("E20").Value = 93

# Row 21
$ws.Range("B21").Value = 0.6392098024506127
$ws.Range("C21").Value = 0.6344086021505376
$ws.Range("D21").Value = 0.6341549540946755
$ws.Range("E21").Value = 93

# Row 22
$ws.Range("B22").Value = 0.6935483870967742
$ws.Range("C22").Value = 0.8775510204081632
$ws.Range("D22").Value = 0.7747747747747746
$ws.Range("E22").Value = 49

# Row 23
$ws.Range("B23").Value = 0.8064516129032258
$ws.Range("C23").Value = 0.5681818181818182
$ws.Range("D23").Value = 0.6666666666666667
$ws.Range("E23").Value = 44

# Row 24
$ws.Range("B24").Value = 0.7311827956989247
$ws.Range("C24").Value = 0.7311827956989247
$ws.Range("D24").Value = 0.7311827956989247
$ws.Range("E24").Value = 0.7311827956989247

# Row 25
$ws.Range("B25").Value = 0.75
$ws.Range("C25").Value = 0.7228664192949907
$ws.Range("D25").Value = 0.7207207207207207
$ws.Range("E25").Value = 93

# Row 26
$ws.Range("B26").Value = 0.7469649670482137
$ws.Range("C26").Value = 0.7311827956989247
$ws.Range("D26").Value = 0.7236268526591106
$ws.Range("E26").Value = 93

Write-Output "Updated classification report values for rows 2-26"